$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp value in M2 (Fecha column) to reflect the new run time
$ws.Range("M2").Value = "2 jul. 2023, 17:18:17"

# Update the selected/active cell in the sheet view
$ws.Range("E3").Select()
